# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45171 (2023-09-02) to 45172 (2023-09-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Data rows span from row 2 to row 533 (row 1 is the header row).
$ws.Range("C2:C533").Value = 45172
